$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column P holds shipment tracking numbers stored as text (not numbers).
# A bare numeric-looking value would be auto-converted to a Number by
# Excel's normal "smart" input parsing, so we force text entry with a
# leading apostrophe and then restore the "Normal" style so the cell
# doesn't pick up a lingering quote-prefix / text number-format style.
$ws.Range("P2").Value = "'320018099707"
$ws.Range("P2").Style = "Normal"

$ws.Range("P3").Value = "'320018114229"
$ws.Range("P3").Style = "Normal"

$ws.Range("P4").Value = "'320018110017"
$ws.Range("P4").Style = "Normal"
